$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "27.899.71"
$ws.Cells.Item(2, 5).Value = "  +4.92%  "
$ws.Cells.Item(3, 4).Value = "1.775.26"
$ws.Cells.Item(3, 5).Value = "  +3.36%  "
$ws.Cells.Item(4, 5).Value = "  +0.21%  "
$ws.Cells.Item(5, 4).Value = "243.09"
$ws.Cells.Item(5, 5).Value = "  +1.08%  "
$ws.Cells.Item(6, 4).Value = "1.001"
$ws.Cells.Item(6, 5).Value = "  +0.20%  "
$ws.Cells.Item(7, 4).Value = "0.4883"
$ws.Cells.Item(7, 5).Value = "  -0.69%  "
$ws.Cells.Item(8, 4).Value = "0.2651"
$ws.Cells.Item(8, 5).Value = "  +2.01%  "
$ws.Cells.Item(9, 5).Value = "  +0.38%  "
$ws.Cells.Item(10, 4).Value = "1.777.73"
$ws.Cells.Item(10, 5).Value = "  +3.43%  "
$ws.Cells.Item(11, 4).Value = "16.27"
$ws.Cells.Item(11, 5).Value = "  +3.61%  "
$ws.Cells.Item(12, 4).Value = "0.07009"
$ws.Cells.Item(12, 5).Value = "  +0.12%  "
$ws.Cells.Item(13, 4).Value = "0.6158"
$ws.Cells.Item(13, 5).Value = "  +1.45%  "
$ws.Cells.Item(14, 4).Value = "4.604"
$ws.Cells.Item(14, 5).Value = "  +2.88%  "
$ws.Cells.Item(15, 4).Value = "79.39"
$ws.Cells.Item(15, 5).Value = "  +3.45%  "
$ws.Cells.Item(16, 4).Value = "27.885.89"
$ws.Cells.Item(16, 5).Value = "  +5.43%  "
$ws.Cells.Item(17, 4).Value = "1.001"
$ws.Cells.Item(17, 5).Value = "  +0.19%  "
$ws.Cells.Item(18, 4).Value = "'1.000"
$ws.Cells.Item(18, 5).Value = "  +0.16%  "
$ws.Cells.Item(19, 5).Value = "  +0.66%  "
$ws.Cells.Item(20, 5).Value = "  +3.74%  "
$ws.Cells.Item(21, 4).Value = "2.008.36"
$ws.Cells.Item(21, 5).Value = "  +2.98%  "
$ws.Cells.Item(22, 4).Value = "4.558"
$ws.Cells.Item(22, 5).Value = "  +3.48%  "
$ws.Cells.Item(23, 4).Value = "8.624"
$ws.Cells.Item(23, 5).Value = "  +1.65%  "
$ws.Cells.Item(24, 4).Value = "5.182"
$ws.Cells.Item(24, 5).Value = "  +1.88%  "
$ws.Cells.Item(25, 4).Value = "141.89"
$ws.Cells.Item(25, 5).Value = "  +3.08%  "
$ws.Cells.Item(26, 4).Value = "15.56"
$ws.Cells.Item(26, 5).Value = "  +1.95%  "
$ws.Cells.Item(27, 4).Value = "1.857"
$ws.Cells.Item(27, 5).Value = "  +6.63%  "
$ws.Cells.Item(28, 2).Value = "BitcoinCash"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(28, 4).Value = "108.92"
$ws.Cells.Item(28, 5).Value = "  +3.01%  "
$ws.Cells.Item(29, 2).Value = "Toncoin"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(29, 4).Value = "1.393"
$ws.Cells.Item(29, 5).Value = "  -0.59%  "
$ws.Cells.Item(30, 4).Value = "4.087"
$ws.Cells.Item(30, 5).Value = "  +4.47%  "
$ws.Cells.Item(31, 4).Value = "0.08231"
$ws.Cells.Item(31, 5).Value = "  +3.68%  "
$ws.Cells.Item(32, 4).Value = "3.762"
$ws.Cells.Item(32, 5).Value = "  +3.44%  "
$ws.Cells.Item(33, 4).Value = "'0.04740"
$ws.Cells.Item(33, 5).Value = "  +5.31%  "
$ws.Cells.Item(34, 4).Value = "1.053"
$ws.Cells.Item(34, 5).Value = "  +5.59%  "
$ws.Cells.Item(35, 4).Value = "2.595"
$ws.Cells.Item(35, 5).Value = "  -0.77%  "
$ws.Cells.Item(36, 4).Value = "0.6393"
$ws.Cells.Item(36, 5).Value = "  +2.27%  "
$ws.Cells.Item(37, 4).Value = "0.9402"
$ws.Cells.Item(37, 5).Value = "  +0.05%  "
$ws.Cells.Item(38, 4).Value = "2.584"
$ws.Cells.Item(39, 4).Value = "2.044"
$ws.Cells.Item(39, 5).Value = "  +1.69%  "
$ws.Cells.Item(40, 4).Value = "'5.870"
$ws.Cells.Item(40, 5).Value = "  +6.41%  "
$ws.Cells.Item(41, 4).Value = "0.01533"
$ws.Cells.Item(41, 5).Value = "  +2.19%  "
$ws.Cells.Item(42, 5).Value = "  +0.25%  "
$ws.Cells.Item(43, 4).Value = "100.25"
$ws.Cells.Item(43, 5).Value = "  +0.71%  "
$ws.Cells.Item(44, 4).Value = "0.3932"
$ws.Cells.Item(44, 5).Value = "  +2.47%  "
$ws.Cells.Item(45, 4).Value = "7.151"
$ws.Cells.Item(45, 5).Value = "  +3.11%  "
$ws.Cells.Item(46, 4).Value = "0.1188"
$ws.Cells.Item(46, 5).Value = "  +3.06%  "
$ws.Cells.Item(47, 4).Value = "'0.05410"
$ws.Cells.Item(47, 5).Value = "  +0.73%  "
$ws.Cells.Item(48, 4).Value = "7.891"
$ws.Cells.Item(48, 5).Value = "  +1.19%  "
$ws.Cells.Item(49, 4).Value = "30.37"
$ws.Cells.Item(49, 5).Value = "  +0.68%  "
$ws.Cells.Item(50, 4).Value = "1.274"
$ws.Cells.Item(50, 5).Value = "  +4.38%  "
$ws.Cells.Item(51, 4).Value = "52.36"
$ws.Cells.Item(51, 5).Value = "  +1.87%  "
